$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91 (shifts existing rows 91-169 down to 92-170)
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new weekly record
$ws.Range("A91").Value = 3
$ws.Range("B91").Value = "Femacal de La Calera"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 44658
$ws.Range("E91").Value = 5
$ws.Range("F91").Value = 100112052
$ws.Range("G91").Value = "Albahaca"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 75
$ws.Range("K91").Value = 4000
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = 4000
$ws.Range("N91").Value = "`$/docena de matas"
$ws.Range("O91").Value = "Provincia de Quillota"
$ws.Range("P91").Value = 667
$ws.Range("Q91").Value = 6
$ws.Range("R91").Value = "Hortaliza"
